$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Update header row (row 1): id/resp1..resp9 -> ID/P1..P9 ---
$headers = @("ID", "P1", "P2", "P3", "P4", "P5", "P6", "P7", "P8", "P9")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Center the header row and give it a plain black font color
$headerRange = $ws.Range("A1:J1")
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4108    # xlCenter
$headerRange.Font.Color = 0

# --- Fill in the previously-empty A2 / E2 cells with "NA" ---
$ws.Cells.Item(2, 1).Value = "NA"
$ws.Cells.Item(2, 5).Value = "NA"

# --- Update the selection to match the header row ---
$headerRange.Select()
